# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45186 (2023-09-17) to serial date 45188 (2023-09-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 171 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cur = $cell.Value2
    if ($cur -ne $null -and [double]$cur -eq 45186) {
        $cell.Value2 = 45188
    }
}
